# Kagiso Rabada.xlsx - append the remaining Delhi Capitals 2020 IPL matches
# to the per-match batting log and replace the single placeholder row with
# the full 7-match data set (rows 2-8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array: dateOfMatch, venueOfMatch, matchResult, ownTeam,
#                   opponentTeam, playerName, runs, balls, numberOf4,
#                   numberOf6, sr
$rows = @(
    @(' Oct 24 2020', ' Abu Dhabi',  'KKR won by 59 runs',
      'Delhi Capitals', 'Kolkata Knight Riders', 'Kagiso Rabada ',
      '9', '10', '1', '0', '90.00'),
    @(' Oct 31 2020', ' Dubai (DSC)', 'Mumbai won by 9 wickets (with 34 balls remaining)',
      'Delhi Capitals', 'Mumbai Indians', 'Kagiso Rabada ',
      '12', '7', '0', '1', '171.42'),
    @(' Nov 10 2020', ' Dubai (DSC)', 'Mumbai won by 5 wickets (with 8 balls remaining)',
      'Delhi Capitals', 'Mumbai Indians', 'Kagiso Rabada ',
      '0', '0', '0', '0', '-'),
    @(' Oct 9 2020', ' Sharjah', 'Capitals won by 46 runs',
      'Delhi Capitals', 'Rajasthan Royals', 'Kagiso Rabada ',
      '2', '3', '0', '0', '66.66'),
    @(' Nov 5 2020', ' Dubai (DSC)', 'Mumbai won by 57 runs',
      'Delhi Capitals', 'Mumbai Indians', 'Kagiso Rabada ',
      '15', '15', '2', '0', '100.00'),
    @(' Sep 20 2020', ' Dubai (DSC)', 'Match tied (Capitals won the one-over eliminator)',
      'Delhi Capitals', 'Kings XI Punjab', 'Kagiso Rabada ',
      '0', '0', '0', '0', '-'),
    @(' Sep 29 2020', ' Abu Dhabi', 'Sunrisers won by 15 runs',
      'Delhi Capitals', 'Sunrisers Hyderabad', 'Kagiso Rabada ',
      '15', '7', '1', '1', '214.28')
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]

    # Columns G:K hold numbers-as-text in the source data (e.g. "90.00",
    # "-") - force text formatting first so Excel doesn't re-parse them
    # into real numbers and strip the formatting/leading zeros.
    $ws.Range("G$r`:K$r").NumberFormat = "@"

    $ws.Cells.Item($r, 1).Value = $values[0]
    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
    $ws.Cells.Item($r, 5).Value = $values[4]
    $ws.Cells.Item($r, 6).Value = $values[5]
    $ws.Cells.Item($r, 7).Value = $values[6]
    $ws.Cells.Item($r, 8).Value = $values[7]
    $ws.Cells.Item($r, 9).Value = $values[8]
    $ws.Cells.Item($r, 10).Value = $values[9]
    $ws.Cells.Item($r, 11).Value = $values[10]
}
